$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Remove the Taxon_Local_ID / ${iNaturalistTaxonId} column entirely.
$ws.Range("A1").EntireColumn.Delete()

# Remove the suborder / infraorder / superfamily columns entirely.
# (After the deletion above these three columns now sit at AQ:AS.)
$ws.Range("AQ1:AS1").EntireColumn.Delete()

# Fix the ${summary.Author} placeholder -> ${summary.authority}.
# (After both column deletions this cell now sits at AX2.)
$ws.Range("AX2").Value = "`${summary.authority}"
